$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SJSU"
$ws.Range("A3").Value = "MLK Jr Library"

$ws.Range("A3").Select()
